# Fruta / hortaliza, semanal
# Permute the weekly data rows (2-6): the D,M,N,O,P,Q,S,T values of each
# row are replaced by those that used to belong to another row, per the
# mapping: target row <- source row
#   2 <- 6
#   3 <- 4
#   4 <- 5
#   5 <- 3
#   6 <- 2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for the columns that move, for rows 2..6.
$cols = @("D","M","N","O","P","Q","S","T")
$original = @{}
foreach ($r in 2..6) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$r").Value2
    }
    $original[$r] = $rowData
}

# Mapping of destination row -> source row (data that should end up there)
$mapping = @{
    2 = 6
    3 = 4
    4 = 5
    5 = 3
    6 = 2
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcData = $original[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $srcData[$col]
    }
}
